$wb = $excel.ActiveWorkbook

# Sheet: Summary
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = 43648.21912747356
$ws.Range("B7").Value = 10346443.45583962
$ws.Range("B8").Value = 24622076.08006534
$ws.Range("B10").Value = 2989385.267199143

# Sheet: Fed-in Capacity
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("M18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("Q20").Value = 150.3839754851235
$ws.Range("M20").Value = 113.4004983079896
$ws.Range("R20").Value = 65.71641987298243
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 117.8828208804077
$ws.Range("J21").Value = 93.17061249236157
$ws.Range("P21").Value = 65.92768427608706
$ws.Range("N21").Value = 38.66169381481656
$ws.Range("N22").Value = 81.96869489115805
$ws.Range("K22").Value = 94.30397654773019
$ws.Range("L22").Value = 90.4687457914608
$ws.Range("N24").Value = 0
$ws.Range("Q24").Value = 94.49434172313325
$ws.Range("P24").Value = 65.92768427608706
$ws.Range("N26").Value = 110.5750244233121
$ws.Range("M26").Value = 113.4004983079896
$ws.Range("K26").Value = 0
$ws.Range("J26").Value = 124.5190384721106
$ws.Range("O26").Value = 117.8828208804077
$ws.Range("I27").Value = 10.12574714858493
$ws.Range("K27").Value = 80.29914934735042
$ws.Range("O27").Value = 0
$ws.Range("Q27").Value = 94.49434172313325
$ws.Range("P27").Value = 0
$ws.Range("K28").Value = 94.30397654773019
$ws.Range("P28").Value = 101.5955875616828
$ws.Range("P29").Value = 135.4597561231036
$ws.Range("O29").Value = 0
$ws.Range("N30").Value = 38.66169381481656
$ws.Range("M30").Value = 51.84373129681028
$ws.Range("M32").Value = 113.4004983079896
$ws.Range("K32").Value = 135.370731907559
$ws.Range("J32").Value = 124.5190384721106
$ws.Range("P32").Value = 135.4597561231036
$ws.Range("Q32").Value = 150.3839754851235
$ws.Range("K33").Value = 80.29914934735042
$ws.Range("L33").Value = 61.18167021676314
$ws.Range("R33").Value = 45.52166981132082
$ws.Range("P33").Value = 65.92768427608706
$ws.Range("O33").Value = 57.81213424001893
$ws.Range("N34").Value = 81.96869489115805
$ws.Range("K34").Value = 94.30397654773019
$ws.Range("M34").Value = 92.09541281912071
$ws.Range("L34").Value = 90.4687457914608
$ws.Range("O34").Value = 96.22962838366004
$ws.Range("K35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("R44").Value = 65.71641987298243
$ws.Range("K46").Value = 94.30397654773019

# Sheet: Unmet Demand
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("M18").Value = 51.84373129681028
$ws.Range("L18").Value = 61.18167021676314
$ws.Range("P18").Value = 65.92768427608706
$ws.Range("O19").Value = 96.22962838366004
$ws.Range("K19").Value = 94.30397654773019
$ws.Range("M19").Value = 92.09541281912071
$ws.Range("L19").Value = 90.4687457914608
$ws.Range("Q20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("R20").Value = 108.0327934026353
$ws.Range("N20").Value = 110.5750244233121
$ws.Range("O20").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N24").Value = 38.66169381481656
$ws.Range("Q24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("K26").Value = 135.370731907559
$ws.Range("J26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("I27").Value = 77.12765456497084
$ws.Range("K27").Value = 0
$ws.Range("O27").Value = 57.81213424001893
$ws.Range("Q27").Value = 0
$ws.Range("P27").Value = 65.92768427608706
$ws.Range("K28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("O29").Value = 117.8828208804077
$ws.Range("N30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("R33").Value = 78.03303713061706
$ws.Range("P33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("K35").Value = 135.370731907559
$ws.Range("Q35").Value = 150.3839754851235
$ws.Range("P35").Value = 135.4597561231036
$ws.Range("L36").Value = 61.18167021676314
$ws.Range("P39").Value = 65.92768427608706
$ws.Range("R44").Value = 108.0327934026353
$ws.Range("K46").Value = 0

# Sheet: Household Surplus
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B7").Value = 328154.0515658425
$ws.Range("B8").Value = 341699.3756980308
$ws.Range("B9").Value = 157065.0487020995
$ws.Range("B10").Value = 303453.8470174211
$ws.Range("B11").Value = 309481.609335087
$ws.Range("B12").Value = 300361.1817230916
$ws.Range("B13").Value = 217059.9772503463
$ws.Range("B14").Value = 158131.4114103872
$ws.Range("B16").Value = 166700.9397412791

# Sheet: Costs and Revenues
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("K2").Value = 85607.3981629999
$ws.Range("P2").Value = 50932.09269021799
$ws.Range("G2").Value = 90142.13413332625
$ws.Range("N2").Value = 48850.92152414424
$ws.Range("H2").Value = 93431.7128511434
$ws.Range("J2").Value = 84143.51302870964
$ws.Range("M2").Value = 63162.14465670575
$ws.Range("L2").Value = 83392.4371715153
$ws.Range("I2").Value = 48591.94772356008
$ws.Range("E3").Value = 133100.0000000001
$ws.Range("D4").Value = 8982.957139551894
$ws.Range("M4").Value = 20142.70273238149
$ws.Range("L4").Value = 40372.99524719104
$ws.Range("I4").Value = 5572.505799235807
$ws.Range("K4").Value = 42587.95623867566
$ws.Range("P4").Value = 7912.650765893709
$ws.Range("G4").Value = 47122.69220900199
$ws.Range("N4").Value = 5831.479599819972
$ws.Range("H4").Value = 50412.27092681914
$ws.Range("J4").Value = 41124.07110438539
$ws.Range("E6").Value = -93692.17018068412
$ws.Range("H6").Value = 39407.82981931594
$ws.Range("K6").Value = 39407.82981931592
$ws.Range("J6").Value = 39407.82981931593
$ws.Range("L6").Value = 39407.82981931594
$ws.Range("P6").Value = 39407.82981931595
